$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update GAME NAME / GAME NUMBER / TOP PRIZES REMAINING cells
$ws.Range('C5').Value = 'Double Pay Day'
$ws.Range("D5").Value = 1062
$ws.Range("E5").Value = 3
$ws.Range('C6').Value = 'Lucky 7'
$ws.Range("D6").Value = 1051
$ws.Range("E6").Value = 5
$ws.Range('C7').Value = '#scratch'
$ws.Range("D7").Value = 1032
$ws.Range('C8').Value = 'Win it All'
$ws.Range("D8").Value = 984
$ws.Range('C9').Value = 'Holiday Cheer'
$ws.Range("D9").Value = 1070
$ws.Range('C10').Value = '$$$'
$ws.Range("D10").Value = 1085
$ws.Range("E10").Value = 4
$ws.Range('C11').Value = 'Money Bags'
$ws.Range("D11").Value = 1007
$ws.Range("E11").Value = 2
$ws.Range('C15').Value = '3 Times Lucky'
$ws.Range("D15").Value = 1116
$ws.Range("E15").Value = 14
$ws.Range('C16').Value = 'Neon 1s'
$ws.Range("D16").Value = 1093
$ws.Range("E16").Value = 3
$ws.Range('C17').Value = '3-2-Won!'
$ws.Range("D17").Value = 1097
$ws.Range("E17").Value = 4
$ws.Range('C26').Value = 'Wild Numbers 10X'
$ws.Range("D26").Value = 1019
$ws.Range("E26").Value = 1
$ws.Range('C27').Value = '4 Way Cash'
$ws.Range("D27").Value = 976
$ws.Range("E27").Value = 2
$ws.Range('C31').Value = 'Jumbo Bucks'
$ws.Range("D31").Value = 1057
$ws.Range('C32').Value = '$20,000 Wild Cherries'
$ws.Range("D32").Value = 1075
$ws.Range('C33').Value = '$100 Loaded'
$ws.Range("D33").Value = 1094
$ws.Range('C35').Value = 'Joker''s Wild'
$ws.Range("D35").Value = 1042
$ws.Range('C37').Value = 'Double Match'
$ws.Range("D37").Value = 1063
$ws.Range('C38').Value = 'Holiday Spectacular'
$ws.Range("D38").Value = 1071
$ws.Range('C65').Value = 'Wild Win!'
$ws.Range("D65").Value = 1012
$ws.Range("E65").Value = 2
$ws.Range('C66').Value = 'Ca$h In'
$ws.Range("D66").Value = 1064
$ws.Range('C67').Value = 'Triple Play'
$ws.Range("D67").Value = 1049
$ws.Range("E67").Value = 1
$ws.Range('C78').Value = '$50 & $100 Blowout'
$ws.Range("D78").Value = 1100
$ws.Range("E78").Value = 127773
$ws.Range('C79').Value = '$500 Fully Loaded'
$ws.Range("D79").Value = 1102
$ws.Range("E79").Value = 2150
$ws.Range("E82").Value = 105

# Update LAST SCRAPE DATE cells (kept as text, matching original shared-string storage)
$dateCells = @('F2','F3','F4','F5','F6','F7','F8','F9','F10','F11','F12','F13','F14','F15','F16','F17','F18','F19','F20','F21','F22','F23','F24','F26','F27','F28','F29','F30','F31','F32','F33','F34','F35','F36','F37','F38','F39','F40','F41','F43','F44','F46','F47','F48','F49','F50','F52','F53','F54','F55','F56','F57','F58','F59','F60','F61','F62','F63','F64','F65','F66','F67','F68','F69','F70','F71','F72','F73','F74','F75','F76','F77','F78','F79','F80','F81','F82','F83','F84','F85','F86','F87','F88','F89')
foreach ($cell in $dateCells) {
    $ws.Range($cell).NumberFormat = "@"
}
$ws.Range('F2').Value = '2019-03-12'
$ws.Range('F3').Value = '2019-03-12'
$ws.Range('F4').Value = '2019-03-12'
$ws.Range('F5').Value = '2019-03-12'
$ws.Range('F6').Value = '2019-03-12'
$ws.Range('F7').Value = '2019-03-12'
$ws.Range('F8').Value = '2019-03-12'
$ws.Range('F9').Value = '2019-03-12'
$ws.Range('F10').Value = '2019-03-12'
$ws.Range('F11').Value = '2019-03-12'
$ws.Range('F12').Value = '2019-03-12'
$ws.Range('F13').Value = '2019-03-12'
$ws.Range('F14').Value = '2019-03-12'
$ws.Range('F15').Value = '2019-03-12'
$ws.Range('F16').Value = '2019-03-12'
$ws.Range('F17').Value = '2019-03-12'
$ws.Range('F18').Value = '2019-03-12'
$ws.Range('F19').Value = '2019-03-12'
$ws.Range('F20').Value = '2019-03-12'
$ws.Range('F21').Value = '2019-03-12'
$ws.Range('F22').Value = '2019-03-12'
$ws.Range('F23').Value = '2019-03-12'
$ws.Range('F24').Value = '2019-03-12'
$ws.Range('F26').Value = '2019-03-12'
$ws.Range('F27').Value = '2019-02-19'
$ws.Range('F28').Value = '2019-03-12'
$ws.Range('F29').Value = '2019-03-12'
$ws.Range('F30').Value = '2019-03-12'
$ws.Range('F31').Value = '2019-03-12'
$ws.Range('F32').Value = '2019-03-12'
$ws.Range('F33').Value = '2019-03-12'
$ws.Range('F34').Value = '2019-03-12'
$ws.Range('F35').Value = '2019-03-05'
$ws.Range('F36').Value = '2019-03-12'
$ws.Range('F37').Value = '2019-03-12'
$ws.Range('F38').Value = '2019-03-12'
$ws.Range('F39').Value = '2019-03-12'
$ws.Range('F40').Value = '2019-03-12'
$ws.Range('F41').Value = '2019-03-12'
$ws.Range('F43').Value = '2019-03-12'
$ws.Range('F44').Value = '2019-03-12'
$ws.Range('F46').Value = '2019-03-12'
$ws.Range('F47').Value = '2019-03-12'
$ws.Range('F48').Value = '2019-03-12'
$ws.Range('F49').Value = '2019-03-12'
$ws.Range('F50').Value = '2019-03-12'
$ws.Range('F52').Value = '2019-03-12'
$ws.Range('F53').Value = '2019-03-12'
$ws.Range('F54').Value = '2019-03-12'
$ws.Range('F55').Value = '2019-03-12'
$ws.Range('F56').Value = '2019-03-12'
$ws.Range('F57').Value = '2019-03-12'
$ws.Range('F58').Value = '2019-03-12'
$ws.Range('F59').Value = '2019-03-12'
$ws.Range('F60').Value = '2019-03-12'
$ws.Range('F61').Value = '2019-03-12'
$ws.Range('F62').Value = '2019-03-12'
$ws.Range('F63').Value = '2019-03-12'
$ws.Range('F64').Value = '2019-03-12'
$ws.Range('F65').Value = '2019-03-12'
$ws.Range('F66').Value = '2019-03-12'
$ws.Range('F67').Value = '2019-03-12'
$ws.Range('F68').Value = '2019-03-12'
$ws.Range('F69').Value = '2019-03-12'
$ws.Range('F70').Value = '2019-03-12'
$ws.Range('F71').Value = '2019-03-12'
$ws.Range('F72').Value = '2019-03-12'
$ws.Range('F73').Value = '2019-03-12'
$ws.Range('F74').Value = '2019-03-12'
$ws.Range('F75').Value = '2019-03-12'
$ws.Range('F76').Value = '2019-03-12'
$ws.Range('F77').Value = '2019-03-12'
$ws.Range('F78').Value = '2019-03-12'
$ws.Range('F79').Value = '2019-03-12'
$ws.Range('F80').Value = '2019-03-12'
$ws.Range('F81').Value = '2019-03-12'
$ws.Range('F82').Value = '2019-03-12'
$ws.Range('F83').Value = '2019-03-12'
$ws.Range('F84').Value = '2019-03-12'
$ws.Range('F85').Value = '2019-03-12'
$ws.Range('F86').Value = '2019-03-12'
$ws.Range('F87').Value = '2019-03-12'
$ws.Range('F88').Value = '2019-03-12'
$ws.Range('F89').Value = '2019-03-12'
foreach ($cell in $dateCells) {
    $ws.Range($cell).ClearFormats()
}
